$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "427. Construct Quad Tree"
$ws.Range("B35").Value = "Medium"
$ws.Range("C35").Value = "Divide and Conquer"
$ws.Range("D35").Value = "Straightforward OOP and grid recursion. Have a helper function and use a 2d loop to check if all values are the same. Use the same input grid and pointers to recursively call each subgrid."
$ws.Range("E35").Value = "https://leetcode.com/problems/construct-quad-tree/solutions/3234703/clean-codes-full-explanation-helper-method-c-java-python3/?envType=study-plan-v2&envId=top-interview-150 "

$ws.Range("B35").Interior.Color = 49407

$ws.Hyperlinks.Add($ws.Range("E35"), "https://leetcode.com/problems/construct-quad-tree/solutions/3234703/clean-codes-full-explanation-helper-method-c-java-python3/?envType=study-plan-v2&envId=top-interview-150 ")
$ws.Range("E35").Style = "Hyperlink"

Write-Host "done"
